$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 16
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("F2").Value = "2025-03-23 21:55:01"

# Update row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = "2025-03-23 21:59:22"

# Delete rows 4-8 (remove entire rows)
$ws.Range("A4:G8").EntireRow.Delete()
